$wb = $excel.ActiveWorkbook

# --- Rename header cells on the existing sheets ---
$wsWeekly = $wb.Worksheets.Item(1)
$wsWeekly.Range("B1").Value = "Weekly_PO_Qty"

$wsMonthly = $wb.Worksheets.Item(2)
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# --- Add the new "PO Forecast" sheet after the last existing sheet ---
# Duplicate the Monthly Trend sheet so the new sheet inherits the same
# sheet-level formatting (margins, outline props, date/header styles)
# instead of starting from a blank engine-default sheet.
$wsMonthly.Copy([System.Reflection.Missing]::Value, $wsMonthly)
$wsForecast = $wb.Worksheets.Item(3)
$wsForecast.Name = "PO Forecast"

# Extend the bold/bordered header styling to columns C and D
$wsForecast.Range("A1:B1").Copy($wsForecast.Range("C1:D1"))

# Extend the date-formatted style down through row 15
$wsForecast.Range("A2").Copy($wsForecast.Range("A2:A15"))

$wsForecast.Range("A1").Value = "ds"
$wsForecast.Range("B1").Value = "PO_Forecast"
$wsForecast.Range("C1").Value = "yhat_lower"
$wsForecast.Range("D1").Value = "yhat_upper"

$rows = @(
    @(45186.99999999999, 22, -30.46048657710787, 75.07749855896402),
    @(45193.99999999999, 23, -31.16888998044095, 79.5716654026804),
    @(45200.99999999999, 24, -25.38601782918717, 76.8580363959009),
    @(45445.99999999999, 61, 7.937462692982734, 116.7785862250356),
    @(45536.99999999999, 75, 20.50289396140212, 129.3053753461474),
    @(45641.99999999999, 91, 38.90502809866103, 143.5418668951559),
    @(45648.99999999999, 92, 35.50544832335577, 146.4817285614184),
    @(45655.99999999999, 93, 36.09570170594329, 149.2126528622291),
    @(45662.99999999999, 95, 38.68652676069535, 149.1198974274585),
    @(45669.99999999999, 96, 43.62882901534231, 147.0568102361585),
    @(45676.99999999999, 97, 44.23574883105031, 150.6058227073527),
    @(45683.99999999999, 98, 40.69792297432461, 152.3243109943162),
    @(45690.99999999999, 99, 43.23683006053418, 152.0628891693325),
    @(45697.99999999999, 100, 44.97036086358842, 153.6812268412772)
)

for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $i + 2
    $row = $rows[$i]
    $wsForecast.Cells.Item($r, 1).Value = $row[0]
    $wsForecast.Cells.Item($r, 2).Value = $row[1]
    $wsForecast.Cells.Item($r, 3).Value = $row[2]
    $wsForecast.Cells.Item($r, 4).Value = $row[3]
}
